$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.197497367858887
$ws.Range("B1").Value = 2.298686981201172
$ws.Range("C1").Value = 6.312384128570557
$ws.Range("D1").Value = 2.03056812286377
$ws.Range("E1").Value = 1.180909156799316
